$wb = $excel.ActiveWorkbook

# Update "Last Updated" timestamp on the Metadata sheet
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 10:09 AM"

# Update the "1 Year" performance column (F) on the Industry Analysis sheet
$ws = $wb.Worksheets.Item("Industry Analysis")
$ws.Cells.Item(2, 6).Value = 21.0016
$ws.Cells.Item(3, 6).Value = -16.2396
$ws.Cells.Item(4, 6).Value = 27.1317
$ws.Cells.Item(5, 6).Value = -50.6494
$ws.Cells.Item(6, 6).Value = 53.2813
$ws.Cells.Item(7, 6).Value = -8.106199999999999
$ws.Cells.Item(8, 6).Value = -9.552099999999999
$ws.Cells.Item(9, 6).Value = 36.3756
$ws.Cells.Item(10, 6).Value = -6.1314
$ws.Cells.Item(11, 6).Value = 31.9081
$ws.Cells.Item(12, 6).Value = -18.4955
$ws.Cells.Item(13, 6).Value = 14.0155
$ws.Cells.Item(14, 6).Value = -36.0718
$ws.Cells.Item(15, 6).Value = -0.1622
$ws.Cells.Item(16, 6).Value = 0.1459
$ws.Cells.Item(17, 6).Value = -22.0012
$ws.Cells.Item(18, 6).Value = 1.0561
$ws.Cells.Item(19, 6).Value = -27.708
$ws.Cells.Item(20, 6).Value = 47.7309
$ws.Cells.Item(21, 6).Value = 12.0959
$ws.Cells.Item(22, 6).Value = 95.1491
$ws.Cells.Item(23, 6).Value = -50.2657
$ws.Cells.Item(24, 6).Value = -13.3427
$ws.Cells.Item(25, 6).Value = -9.9316
$ws.Cells.Item(26, 6).Value = 5.8244
$ws.Cells.Item(27, 6).Value = -32.7692
$ws.Cells.Item(28, 6).Value = -24.8224
$ws.Cells.Item(29, 6).Value = -18.4191
$ws.Cells.Item(30, 6).Value = 25.8569
$ws.Cells.Item(31, 6).Value = 58.4712
$ws.Cells.Item(32, 6).Value = -3.3862
$ws.Cells.Item(33, 6).Value = -6.3282
$ws.Cells.Item(34, 6).Value = 27.7203
$ws.Cells.Item(35, 6).Value = 4.4873
$ws.Cells.Item(36, 6).Value = -4.9458
$ws.Cells.Item(37, 6).Value = 3.6074
$ws.Cells.Item(38, 6).Value = -23.3973
$ws.Cells.Item(39, 6).Value = 8.7355
$ws.Cells.Item(40, 6).Value = -5.8541
$ws.Cells.Item(41, 6).Value = -8.3934
$ws.Cells.Item(42, 6).Value = 20.3818
$ws.Cells.Item(43, 6).Value = 14.3164
$ws.Cells.Item(44, 6).Value = -12.6846
$ws.Cells.Item(45, 6).Value = 28.4075
$ws.Cells.Item(46, 6).Value = -1.1135
$ws.Cells.Item(47, 6).Value = -37.1997
$ws.Cells.Item(48, 6).Value = -29.8569
$ws.Cells.Item(49, 6).Value = -27.5511
$ws.Cells.Item(50, 6).Value = -49.7478
$ws.Cells.Item(51, 6).Value = -51.8002
$ws.Cells.Item(52, 6).Value = -38.5254
$ws.Cells.Item(53, 6).Value = -12.4886
$ws.Cells.Item(54, 6).Value = -5.0725
$ws.Cells.Item(55, 6).Value = -17.7445
$ws.Cells.Item(56, 6).Value = -26.636
$ws.Cells.Item(57, 6).Value = -29.3361
$ws.Cells.Item(58, 6).Value = -11.9574
$ws.Cells.Item(59, 6).Value = -24.5687
$ws.Cells.Item(60, 6).Value = -12.3
$ws.Cells.Item(61, 6).Value = -10.9446
$ws.Cells.Item(62, 6).Value = -17.1229
$ws.Cells.Item(63, 6).Value = -9.5038
$ws.Cells.Item(64, 6).Value = 54.2749
$ws.Cells.Item(65, 6).Value = -43.4736
$ws.Cells.Item(66, 6).Value = 13.2687
$ws.Cells.Item(67, 6).Value = 12.7149
$ws.Cells.Item(68, 6).Value = 24.8057
$ws.Cells.Item(69, 6).Value = -17.0328
$ws.Cells.Item(70, 6).Value = -6.8927
$ws.Cells.Item(71, 6).Value = 13.6034
$ws.Cells.Item(72, 6).Value = 3.9995
$ws.Cells.Item(73, 6).Value = -16.226
$ws.Cells.Item(74, 6).Value = -16.2448
$ws.Cells.Item(75, 6).Value = 28.6924
$ws.Cells.Item(76, 6).Value = 48.9752
